# Renamed Trade type values
#
# The Sample Logbook workbook had a "Trade Type" column using values such as
# "Deposit (INR)", "Withdraw (BTC)", "Reinvest (USDT)", etc., together with a
# 5th sample row per sheet. The sheets are simplified down to three sample
# trade rows using the generic Trade Type values "Buy" / "Reinvest" / "Sell",
# and the Coin column now cycles through ABC / DEF / XYZ (dropping the old
# "IJK" sample and the redundant second "XYZ" row).

$wb = $excel.ActiveWorkbook

$wsUSDT = $wb.Worksheets.Item("USDT")
$wsBTC  = $wb.Worksheets.Item("BTC")

foreach ($ws in @($wsUSDT, $wsBTC)) {
    # Update the three retained sample rows with the new Coin / Trade Type text.
    $ws.Range("B2").Value2 = "ABC"
    $ws.Range("L2").Value2 = "Buy"

    $ws.Range("B3").Value2 = "DEF"
    $ws.Range("L3").Value2 = "Reinvest"

    $ws.Range("B4").Value2 = "XYZ"
    $ws.Range("L4").Value2 = "Sell"

    # Drop the old 5th sample row entirely (shrinks dimension to A1:L4).
    $ws.Rows.Item(5).Delete()
}

# Restore each sheet's own remembered selection.
[void]$wsUSDT.Range("B4").Select()
[void]$wsBTC.Range("E5").Select()

# BTC is now the active/selected tab (was USDT before).
$wsBTC.Activate()
